# Update New Orleans xlsx: add a "State" column to hotel_info, and move
# review_info so it is the first sheet tab (before hotel_info).

$wb = $excel.ActiveWorkbook

$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new "State" column between "Hotel_Name" (col B) and "City" (col C)
$hotelSheet.Columns.Item(3).Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Reorder the sheet tabs so review_info comes before hotel_info
$reviewSheet.Move($hotelSheet)
